$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Price (D) column cells being updated, to preserve exact
# string formatting (trailing zeros, multi-dot thousands separators, etc.)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply cell value updates
$ws.Range("D2").Value = "50.056.16"
$ws.Range("E2").Value = "  +3.68%  "
$ws.Range("D3").Value = "2.648.94"
$ws.Range("E3").Value = "  +5.96%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "114.07"
$ws.Range("E5").Value = "  +7.46%  "
$ws.Range("D6").Value = "326.78"
$ws.Range("E6").Value = "  +2.76%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +3.50%  "
$ws.Range("D10").Value = "40.95"
$ws.Range("E10").Value = "  +5.36%  "
$ws.Range("D11").Value = "20.16"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("E12").Value = "  +2.34%  "
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("D14").Value = "7.40"
$ws.Range("E14").Value = "  +3.96%  "
$ws.Range("D15").Value = "3.063.54"
$ws.Range("E15").Value = "  +5.96%  "
$ws.Range("D16").Value = "2.647.87"
$ws.Range("E16").Value = "  +5.93%  "
$ws.Range("D17").Value = "0.873"
$ws.Range("E17").Value = "  +5.15%  "
$ws.Range("D18").Value = "49.982.07"
$ws.Range("E18").Value = "  +3.89%  "
$ws.Range("D19").Value = "13.24"
$ws.Range("E20").Value = "  +2.85%  "
$ws.Range("D21").Value = "2.92"
$ws.Range("E21").Value = "  -3.76%  "
$ws.Range("D23").Value = "72.28"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").Value = "276.18"
$ws.Range("E24").Value = "  +2.34%  "
$ws.Range("D25").Value = "2.59"
$ws.Range("E25").Value = "  +3.04%  "
$ws.Range("D26").Value = "26.85"
$ws.Range("D27").Value = "1.00"
$ws.Range("E28").Value = "  +3.38%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "36.42"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").Value = "50.14"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("D33").Value = "5.47"
$ws.Range("D34").Value = "19.58"
$ws.Range("E34").Value = "  +2.41%  "
$ws.Range("D35").Value = "0.0816"
$ws.Range("E35").Value = "  +5.44%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").Value = "5.03"
$ws.Range("E37").Value = "  +9.39%  "
$ws.Range("E38").Value = "  +6.74%  "
$ws.Range("E39").Value = "  +7.68%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "0.113"
$ws.Range("E40").Value = "  +1.97%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "123.51"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("D43").Value = "21.96"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("D44").Value = "0.0317"
$ws.Range("E44").Value = "  +4.50%  "
$ws.Range("D45").Value = "2.085.16"
$ws.Range("E45").Value = "  +4.03%  "
$ws.Range("E46").Value = "  +6.30%  "
$ws.Range("D47").Value = "2.28"
$ws.Range("E47").Value = "  +14.34%  "
$ws.Range("D48").Value = "2.00"
$ws.Range("E48").Value = "  +4.46%  "
$ws.Range("D49").Value = "9.16"
$ws.Range("E49").Value = "  +2.65%  "
$ws.Range("D50").Value = "5.39"
$ws.Range("E50").Value = "  +4.25%  "
$ws.Range("D51").Value = "59.80"
$ws.Range("E51").Value = "  +5.40%  "
